# Edit script: add "2022-Q3" data to 06690-海尔智家.xlsx
#
# 1) Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q3 and shift the previously-existing quarter rows down by one.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计",
#    containing the per-fund holdings table for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - shift existing rows down and add 2022-Q3 on top
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# The sheet grows by one row (A1:D7 -> A1:D8). Copy row 7's formatting
# onto the new row 8 first (so the bold/centred index-column style
# carries over), then fill in the data bottom-up. All target values are
# literal (taken straight from the intended final sheet), so the write
# order does not matter for correctness.
$summary.Rows.Item(7).Copy()
$summary.Rows.Item(8).PasteSpecial(-4122)

$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(8,2).Value = "2020-Q4"
$summary.Cells.Item(8,3).Value = 9
$summary.Cells.Item(8,4).Value = 6.99

$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2021-Q1"
$summary.Cells.Item(7,3).Value = 31
$summary.Cells.Item(7,4).Value = 35.31

$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q2"
$summary.Cells.Item(6,3).Value = 15
$summary.Cells.Item(6,4).Value = 22.47

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q3"
$summary.Cells.Item(5,3).Value = 26
$summary.Cells.Item(5,4).Value = 13.58

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2021-Q4"
$summary.Cells.Item(4,3).Value = 27
$summary.Cells.Item(4,4).Value = 15.42

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q1"
$summary.Cells.Item(3,3).Value = 24
$summary.Cells.Item(3,4).Value = 10.97

# New top data row: 2022-Q3. Copy A3's format onto A2 first so the new
# index cell picks up the same bold/centered style as its neighbours.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 39
$summary.Cells.Item(2,4).Value = 15.22

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

$fundData = @(
    ,@(0, '513180', '华夏恒生科技ETF（QDII）', '103.92', '94.41', '4.28', '4.4478', 10)
    ,@(1, '010557', '汇添富数字生活主题六个月持有期混合', '46.22', '89.69', '4.76', '2.2001', 2)
    ,@(2, '011399', '汇添富数字未来混合A', '37.17', '90.07', '4.74', '1.7619', 2)
    ,@(3, '011665', '汇添富数字经济引领发展三年持有期混合A', '65.89', '59.90', '2.51', '1.6538', 5)
    ,@(4, '513010', '易方达恒生科技ETF（QDII）', '25.49', '94.56', '4.28', '1.0910', 10)
    ,@(5, '012588', '南方港股通优势企业混合A', '25.02', '81.50', '3.09', '0.7731', 7)
    ,@(6, '159740', '大成恒生科技ETF（QDII）', '9.95', '93.00', '4.21', '0.4189', 10)
    ,@(7, '012805', '广发恒生科技指数（QDII）C', '7.98', '80.85', '3.66', '0.2921', 10)
    ,@(8, '513580', '华安恒生科技ETF（QDII）', '6.44', '95.77', '4.32', '0.2782', 10)
    ,@(9, '006752', '天弘港股通精选灵活配置混合A', '4.67', '93.60', '5.79', '0.2704', 8)
    ,@(10, '159742', '博时恒生科技ETF（QDII）', '6.13', '95.25', '4.31', '0.2642', 10)
    ,@(11, '011400', '汇添富数字未来混合C', '5.25', '90.07', '4.74', '0.2488', 2)
    ,@(12, '005197', '工银瑞信沪港深精选灵活配置混合A', '4.01', '94.39', '5.23', '0.2097', 6)
    ,@(13, '202801', '南方全球精选配置（QDII-FOF）', '15.80', '29.52', '1.10', '0.1738', 9)
    ,@(14, '006753', '天弘港股通精选灵活配置混合C', '2.52', '93.60', '5.79', '0.1459', 8)
    ,@(15, '005504', '汇添富沪港深大盘价值混合A', '3.14', '91.24', '4.61', '0.1448', 6)
    ,@(16, '159741', '嘉实恒生科技ETF（QDII）', '2.85', '99.55', '4.50', '0.1282', 10)
    ,@(17, '012804', '广发恒生科技指数（QDII）A', '3.18', '80.85', '3.66', '0.1164', 10)
    ,@(18, '513890', '上投摩根恒生科技ETF（QDII）', '1.73', '95.47', '4.83', '0.0836', 9)
    ,@(19, '011666', '汇添富数字经济引领发展三年持有期混合C', '3.02', '59.90', '2.51', '0.0758', 5)
    ,@(20, '007109', '南方沪港深核心优势混合', '1.59', '85.13', '4.11', '0.0653', 6)
    ,@(21, '013127', '汇添富恒生科技指数（QDII）A', '1.49', '91.29', '4.13', '0.0615', 10)
    ,@(22, '013128', '汇添富恒生科技指数（QDII）C', '1.45', '91.29', '4.13', '0.0599', 10)
    ,@(23, '012589', '南方港股通优势企业混合C', '1.75', '81.50', '3.09', '0.0541', 7)
    ,@(24, '005198', '工银瑞信沪港深精选灵活配置混合C', '1.01', '94.39', '5.23', '0.0528', 6)
    ,@(25, '513380', '广发恒生科技（QDII-ETF）', '1.08', '84.71', '3.79', '0.0409', 10)
    ,@(26, '004249', '安信中国制造混合', '0.52', '89.55', '5.02', '0.0261', 7)
    ,@(27, '008254', '华宝致远混合（QDII）C', '0.40', '85.90', '4.14', '0.0166', 8)
    ,@(28, '006205', '汇添富沪港深优势精选定期开放混合', '0.33', '95.18', '4.59', '0.0151', 6)
    ,@(29, '008253', '华宝致远混合（QDII）A', '0.34', '85.90', '4.14', '0.0141', 8)
    ,@(30, '005255', '浦银安盛港股通量化混合A', '0.29', '78.68', '3.41', '0.0099', 10)
    ,@(31, '010777', '浙商智选家居股票A', '0.12', '90.64', '5.52', '0.0066', 7)
    ,@(32, '004321', '前海开源沪港深强国产业灵活配置混合', '0.11', '78.52', '5.87', '0.0065', 1)
    ,@(33, '010778', '浙商智选家居股票C', '0.08', '90.64', '5.52', '0.0044', 7)
    ,@(34, '005707', '富国港股通量化精选股票A', '0.22', '89.83', '1.82', '0.0040', 9)
    ,@(35, '013224', '浦银安盛港股通量化混合C', '0.05', '78.68', '3.41', '0.0017', 10)
    ,@(36, '014163', '富国港股通量化精选股票C', '0.00', '89.83', '1.82', 0, 9)
    ,@(37, '015118', '汇添富沪港深大盘价值混合C', '0.00', '91.24', '4.61', 0, 6)
    ,@(38, '015119', '汇添富沪港深大盘价值混合D', '0.00', '91.24', '4.61', 0, 6)
)

# Columns B-G hold text in the source data (fund codes with leading
# zeros, and numeric-looking percentages/scale figures stored as text),
# so force text formatting before writing those columns. A few rows at
# the bottom have a genuinely-numeric 0 in column G (held market value
# rounds to 0), so those cells are reset back to General afterwards.
$q3.Range("B2:G40").NumberFormat = "@"

foreach ($row in $fundData) {
    $r = [int]$row[0] + 2
    $q3.Cells.Item($r,1).Value = $row[0]
    $q3.Cells.Item($r,2).Value = $row[1]
    $q3.Cells.Item($r,3).Value = $row[2]
    $q3.Cells.Item($r,4).Value = $row[3]
    $q3.Cells.Item($r,5).Value = $row[4]
    $q3.Cells.Item($r,6).Value = $row[5]
    if ($row[6] -is [string]) {
        $q3.Cells.Item($r,7).Value = $row[6]
    } else {
        $q3.Cells.Item($r,7).NumberFormat = "General"
        $q3.Cells.Item($r,7).Value = $row[6]
    }
    $q3.Cells.Item($r,8).Value = $row[7]
}

# Leave the original sheet active/selected, matching the source workbook.
$summary.Activate()
[void]$summary.Range("A1").Select()

Write-Host "Edit applied: 2022-Q3 sheet and summary row inserted."
